$d = $word.ActiveDocument
$emdash = [char]0x2014

# ---------------------------------------------------------------------------
# Step A: Simple one-for-one text replacements (paragraph count unchanged).
# Each Find/Replace is scoped to the specific paragraph's Range so that it
# cannot accidentally match identical text living in another paragraph
# (e.g. "Inject 3 - Automation" also appears in the Title).
# ---------------------------------------------------------------------------

$r = $d.Paragraphs.Item(4).Range
$r.Find.Execute("Jeffrey Fonseca", $false, $false, $false, $false, $false, $true, 1, $false, "Linux Team", 2) | Out-Null

$r = $d.Paragraphs.Item(5).Range
$r.Find.Execute("Linux Team", $false, $false, $false, $false, $false, $true, 1, $false, "Inject 3 $emdash Automation", 2) | Out-Null

$r = $d.Paragraphs.Item(7).Range
$r.Find.Execute("Inject 3 $emdash Automation", $false, $false, $false, $false, $false, $true, 1, $false, "Objective: Create a script or other system that automates updates, ran every 24 hours", 2) | Out-Null

$r = $d.Paragraphs.Item(8).Range
$r.Find.Execute("Objective: Create a script or other system that automates updates, ran every 24 hours", $false, $false, $false, $false, $false, $true, 1, $false, "I have created a systemd timer and unit file that run every 24 hours, designed for a variety of linux distros with different package managers.", 2) | Out-Null

$r = $d.Paragraphs.Item(9).Range
$r.Find.Execute("I have created a systemd timer and unit file that run every 24 hours, designed for a variety of linux distros with different package managers.", $false, $false, $false, $false, $false, $true, 1, $false, "When these files are placed in their respective locations, the system will run the update commands of various pacakge managers. Distros unsupported by a package manager will simply fail to run for that package manager.", 2) | Out-Null

# ---------------------------------------------------------------------------
# Step B: Delete the two paragraphs that are no longer needed, each together
# with its own paragraph mark:
#   - Paragraph 10 still holds the old "...manager.ast" BodyText, whose
#     (corrected) content was already copied forward into paragraph 9 above.
#   - Paragraph 6 still holds the old "8/26/2023" BodyText, which has no
#     counterpart in the new document at all.
# Delete the higher-numbered paragraph first so paragraph 6's index stays
# valid for the second delete.
# ---------------------------------------------------------------------------

$d.Paragraphs.Item(10).Range.Delete() | Out-Null
$d.Paragraphs.Item(6).Range.Delete() | Out-Null

# ---------------------------------------------------------------------------
# Step C: Fix the two small text fragments inside the SourceCode blocks.
# These paragraphs kept their rich (multi-run / styled) content untouched;
# only one run's text changes in each. After the deletions above, the big
# SourceCode block is paragraph 9 and the small one is paragraph 10.
# ---------------------------------------------------------------------------

$r = $d.Paragraphs.Item(9).Range
$r.Find.Execute("Dnf autoupdate service", $false, $false, $false, $false, $false, $true, 1, $false, "autoupdate service", 2) | Out-Null

$r = $d.Paragraphs.Item(10).Range
$r.Find.Execute("Dnf autoupdate timer.", $false, $false, $false, $false, $false, $true, 1, $false, "autoupgrade timer.", 2) | Out-Null
